# Enhance user registration process by adding email, course, class name,
# and roll number fields for new registrants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: sahil / 1122990088 (Name, Mobile only)
$ws.Range("A4").Value = "sahil"
$ws.Range("B4").Value = "'1122990088"

# Row 5: sameer - full registration record (Name, Mobile, Email, Course, Class, Roll No)
$ws.Range("A5").Value = "sameer"
$ws.Range("B5").Value = "'0000000000"
$ws.Range("C5").Value = "mdshaiksahil0510@gmail.com"
$ws.Range("D5").Value = "CSE IOT"
$ws.Range("E5").Value = "CSEIOT"
$ws.Range("F5").Value = "'161023749019"

# Row 6: masood - full registration record (Name, Mobile, Email, Course, Class, Roll No)
$ws.Range("A6").Value = "masood"
$ws.Range("B6").Value = "'7981842202"
$ws.Range("C6").Value = "2005syedmasood@gmail.com"
$ws.Range("D6").Value = "BE"
$ws.Range("E6").Value = "CSE"
$ws.Range("F6").Value = "'161023733094"

# New rows use the workbook default (unstyled) formatting, matching the
# rest of the sheet's plain-registration look.
$ws.Range("A4:B4").Style = "Normal"
$ws.Range("A5:F5").Style = "Normal"
$ws.Range("A6:F6").Style = "Normal"
